# Update the "Förändrad" (Changed) date column (C) for every data row
# (rows 2-295) from 2023-10-03 (serial 45202) to 2023-10-04 (serial 45203).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C2:C295").Value = 45203
